# Atualização de bases das ligas, do dia: 19-06-2024 às 21:51
# Swap the match-data (columns B through AD) between pairs of rows.
# Column A (the running id) stays put; only the rest of each row's
# contents are exchanged with its pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each pair of Excel row numbers whose B:AD contents must be swapped.
$pairs = @(
    @(17, 18),
    @(25, 26),
    @(48, 49),
    @(59, 60),
    @(161, 162)
)

# Columns B (2) through AD (30) inclusive.
$firstCol = 2
$lastCol = 30

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cell1 = $ws.Cells.Item($r1, $col)
        $cell2 = $ws.Cells.Item($r2, $col)

        $v1 = $cell1.Value2
        $v2 = $cell2.Value2

        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}
